$d = $word.ActiveDocument

$replacements = @(
    @{old="543÷4=135, 3"; new="132÷7=18, 6"},
    @{old="475÷4=118, 3"; new="557÷3=185, 2"},
    @{old="765÷2=382, 1"; new="317÷5=63, 2"},
    @{old="558÷3=186, 0"; new="676÷3=225, 1"},
    @{old="567÷7=81, 0"; new="918÷8=114, 6"},
    @{old="854÷5=170, 4"; new="287÷2=143, 1"},
    @{old="816÷2=408, 0"; new="244÷8=30, 4"},
    @{old="834÷3=278, 0"; new="146÷7=20, 6"},
    @{old="262÷2=131, 0"; new="972÷7=138, 6"},
    @{old="753÷6=125, 3"; new="639÷9=71, 0"},
    @{old="635÷7=90, 5"; new="281÷7=40, 1"},
    @{old="465÷7=66, 3"; new="815÷4=203, 3"},
    @{old="882÷5=176, 2"; new="275÷3=91, 2"},
    @{old="457÷3=152, 1"; new="413÷7=59, 0"},
    @{old="108÷3=36, 0"; new="429÷9=47, 6"},
    @{old="506÷3=168, 2"; new="730÷5=146, 0"},
    @{old="925÷3=308, 1"; new="411÷3=137, 0"},
    @{old="273÷7=39, 0"; new="889÷5=177, 4"},
    @{old="647÷5=129, 2"; new="830÷2=415, 0"},
    @{old="158÷3=52, 2"; new="618÷8=77, 2"},
    @{old="295÷4=73, 3"; new="943÷3=314, 1"},
    @{old="347÷3=115, 2"; new="742÷2=371, 0"},
    @{old="147÷9=16, 3"; new="231÷5=46, 1"},
    @{old="240÷4=60, 0"; new="734÷8=91, 6"},
    @{old="276÷2=138, 0"; new="173÷9=19, 2"}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
